function Replace-Range {
    param(
        [object]$Doc,
        [string]$SearchText,
        [string]$XmlFragment
    )
    $full = $Doc.Range(0, $Doc.Content.End)
    $found = $full.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Search text not found: $SearchText"
    }
    $target = $Doc.Range($full.Start, $full.End)
    $target.InsertXML($XmlFragment)
}

function Insert-ParagraphBreakBefore {
    param(
        [object]$Doc,
        [string]$SearchText,
        [string]$PPrXml
    )
    $full = $Doc.Range(0, $Doc.Content.End)
    $found = $full.Find.Execute($SearchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Search text not found: $SearchText"
    }
    $boundary = $full.Start
    $ins = $Doc.Range($boundary, $boundary)
    $xml = "<w:p xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"">$PPrXml</w:p>"
    $ins.InsertXML($xml)
}

$d = $word.ActiveDocument
$W = "xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main"""
$RPR = "<w:rPr><w:rFonts w:ascii=""Verdana"" w:hAnsi=""Verdana"" w:cs=""Verdana""/><w:kern w:val=""0""/></w:rPr>"

# --- Edit 1: insert an empty paragraph before "Le format PDF est un format de fichier qui préserve" ---
$pprXml = "<w:pPr><w:autoSpaceDE w:val=""0""/><w:autoSpaceDN w:val=""0""/><w:adjustRightInd w:val=""0""/><w:spacing w:after=""0"" w:line=""240"" w:lineRule=""auto""/><w:rPr><w:rFonts w:ascii=""Verdana"" w:hAnsi=""Verdana"" w:cs=""Verdana""/><w:kern w:val=""0""/></w:rPr></w:pPr>"
Insert-ParagraphBreakBefore $d "Le format PDF est un format de fichier qui préserve les polices, les images les objets" $pprXml

# --- Edit 2: "... et la plate-forme utilisées pour le créer ..." -> mark "la plate-forme utilisées" as a grammar error ---
$xml2 = "<w:p $W>" + `
        "<w:r>$RPR<w:t xml:space=""preserve""> </w:t></w:r>" + `
        "<w:proofErr w:type=""gramStart""/>" + `
        "<w:r>$RPR<w:t>la plate-forme utilisées</w:t></w:r>" + `
        "<w:proofErr w:type=""gramEnd""/>" + `
        "<w:r>$RPR<w:t xml:space=""preserve""> pour le créer. Les fichiers PDF peuvent être créés avec des</w:t></w:r>" + `
        "</w:p>"
Replace-Range $d " la plate-forme utilisées pour le créer. Les fichiers PDF peuvent être créés avec des" $xml2

# --- Edit 3: "... et systèmes d'exploitations: Mac OS ..." -> mark "d'exploitations:" as a grammar error ---
$xml3 = "<w:p $W>" + `
        "<w:r>$RPR<w:t xml:space=""preserve""> et systèmes </w:t></w:r>" + `
        "<w:proofErr w:type=""gramStart""/>" + `
        "<w:r>$RPR<w:t>d'exploitations:</w:t></w:r>" + `
        "<w:proofErr w:type=""gramEnd""/>" + `
        "<w:r>$RPR<w:t xml:space=""preserve""> Mac OS, Windows, Linux, Palm OS, Pocket PC,</w:t></w:r>" + `
        "</w:p>"
Replace-Range $d " et systèmes d'exploitations: Mac OS, Windows, Linux, Palm OS, Pocket PC," $xml3

# --- Edit 4a: "... dont certains sont des logiciels libres, existent " -> mark "existent ... également" as a grammar error (start) ---
$xml4a = "<w:p $W>" + `
        "<w:r>$RPR<w:t xml:space=""preserve""> dont certains sont des logiciels libres, </w:t></w:r>" + `
        "<w:proofErr w:type=""gramStart""/>" + `
        "<w:r>$RPR<w:t xml:space=""preserve"">existent </w:t></w:r>" + `
        "</w:p>"
Replace-Range $d " dont certains sont des logiciels libres, existent " $xml4a

# --- Edit 4b: "également. La génération de" -> close the grammar error after "également" ---
$xml4b = "<w:p $W>" + `
        "<w:r>$RPR<w:t>également</w:t></w:r>" + `
        "<w:proofErr w:type=""gramEnd""/>" + `
        "<w:r>$RPR<w:t>. La génération de</w:t></w:r>" + `
        "</w:p>"
Replace-Range $d "également. La génération de" $xml4b

Write-Output "All edits applied."
